$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 72 - this shifts the existing rows 72..104 down to 73..105,
# matching the diff (dimension grows from A1:T104 to A1:T105).
$ws.Rows("72:72").Insert()

# Populate the newly inserted row 72 with the new record.
$ws.Range("A72").Value = 4
$ws.Range("B72").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C72").Value = "Los Lagos"
$ws.Range("D72").Value = 44468
$ws.Range("E72").Value = 10
$ws.Range("F72").Value = "Fruta"
$ws.Range("G72").Value = 100102
$ws.Range("H72").Value = "Cítricos"
$ws.Range("I72").Value = 100102004
$ws.Range("J72").Value = "Mandarina"
$ws.Range("K72").Value = "Murcott"
$ws.Range("L72").Value = "Primera"
$ws.Range("M72").Value = 300
$ws.Range("N72").Value = 6500
$ws.Range("O72").Value = 6500
$ws.Range("P72").Value = 6500
$ws.Range("Q72").Value = "$/bandeja 10 kilos"
$ws.Range("R72").Value = "Provincia de Limarí"
$ws.Range("S72").Value = 650
$ws.Range("T72").Value = 10
